$wb = $excel.ActiveWorkbook

# --- Sheet "8 Dena Demas": update the Low/Medium/High expertise ratings
# in column C (Keyword # rows 2-12) ---
$ws3 = $wb.Worksheets.Item("8 Dena Demas")

$ws3.Range("C2").Value = "M"
$ws3.Range("C3").Value = ""
$ws3.Range("C4").Value = "H"
$ws3.Range("C5").Value = "M"
$ws3.Range("C6").Value = "L"
$ws3.Range("C7").Value = "M"
$ws3.Range("C8").Value = "M"
$ws3.Range("C9").Value = "L"
$ws3.Range("C10").Value = ""
$ws3.Range("C11").Value = "M"
$ws3.Range("C12").Value = "L"

# --- Update the on-screen selection state that was left behind on the
# "Conflicts of Interest" and "8 Dena Demas" sheets ---
$ws2 = $wb.Worksheets.Item("Conflicts of Interest")
$ws2.Activate() | Out-Null
$ws2.Range("A10:K10").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("A10:K10").Select() | Out-Null

# Restore the originally active sheet/tab
$ws1 = $wb.Worksheets.Item("Expertise by Keywords - Instr.")
$ws1.Activate() | Out-Null
